$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.398.74'
$ws.Range('E2').Value = '  +4.00%  '

$ws.Range('D3').Value = '2.444.03'
$ws.Range('E3').Value = '  +3.56%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '556.02'
$ws.Range('E5').Value = '  +2.67%  '

$ws.Range('D6').Value = '138.64'
$ws.Range('E6').Value = '  +1.85%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  +1.51%  '

$ws.Range('E9').Value = '  +3.64%  '

$ws.Range('D10').Value = '5.78'
$ws.Range('E10').Value = '  +3.67%  '

$ws.Range('E11').Value = '  +1.29%  '

$ws.Range('E12').Value = '  -2.01%  '

$ws.Range('D13').Value = '''24.90'

$ws.Range('D14').Value = '2.878.67'
$ws.Range('E14').Value = '  +3.56%  '

$ws.Range('D15').Value = '60.316.93'
$ws.Range('E15').Value = '  +3.91%  '

$ws.Range('E16').Value = '  +3.51%  '

$ws.Range('D17').Value = '2.459.52'
$ws.Range('E17').Value = '  +4.03%  '

$ws.Range('D18').Value = '11.39'
$ws.Range('E18').Value = '  +5.90%  '

$ws.Range('D19').Value = '4.42'
$ws.Range('E19').Value = '  +3.17%  '

$ws.Range('D20').Value = '335.32'
$ws.Range('E20').Value = '  +0.59%  '

$ws.Range('D21').Value = '6.92'
$ws.Range('E21').Value = '  +1.68%  '

$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.05%  '

$ws.Range('D23').Value = '64.68'
$ws.Range('E23').Value = '  +2.98%  '

$ws.Range('E24').Value = '  +2.51%  '

$ws.Range('D25').Value = '8.58'
$ws.Range('E25').Value = '  +0.51%  '

$ws.Range('E26').Value = '  +0.06%  '

$ws.Range('E27').Value = '  -0.67%  '

$ws.Range('D28').Value = '0.0₃0792'
$ws.Range('E28').Value = '  +7.06%  '

$ws.Range('D29').Value = '''1.80'
$ws.Range('E29').Value = '  +2.83%  '

$ws.Range('D30').Value = '171.27'
$ws.Range('E30').Value = '  -0.99%  '

$ws.Range('D31').Value = '6.31'
$ws.Range('E31').Value = '  +2.25%  '

$ws.Range('D32').Value = '18.83'
$ws.Range('E32').Value = '  +1.49%  '

$ws.Range('E33').Value = '  -1.26%  '

$ws.Range('E35').Value = '  +5.08%  '

$ws.Range('D36').Value = '4.27'
$ws.Range('E36').Value = '  +0.74%  '

$ws.Range('E37').Value = '  +0.08%  '

$ws.Range('E38').Value = '  +0.35%  '

$ws.Range('D39').Value = '40.09'
$ws.Range('E39').Value = '  +1.77%  '

$ws.Range('D40').Value = '0.415'
$ws.Range('E40').Value = '  +9.35%  '

$ws.Range('D41').Value = '318.14'
$ws.Range('E41').Value = '  +8.16%  '

$ws.Range('D42').Value = '144.06'
$ws.Range('E42').Value = '  -1.25%  '

$ws.Range('E43').Value = '  +1.87%  '

$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '19.94'
$ws.Range('E44').Value = '  +3.49%  '

$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = '0.0963'
$ws.Range('E45').Value = '  +1.40%  '

$ws.Range('E46').Value = '  +4.18%  '

$ws.Range('B47').Value = 'Polygon'
$ws.Range('C47').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D47').Value = '0.411'
$ws.Range('E47').Value = '  +5.41%  '

$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.574'
$ws.Range('E48').Value = '  +1.59%  '

$ws.Range('E49').Value = '  +1.90%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '18.15'
$ws.Range('E50').Value = '  +3.59%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '11.05'
$ws.Range('E51').Value = '  -0.16%  '
